$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44172
$ws.Range("L2").Value = "Especial"
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 6500
$ws.Range("O2").Value = 7000
$ws.Range("P2").Value = 6750
$ws.Range("S2").Value = 2250

$ws.Range("D3").Value = 44172
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 160
$ws.Range("N3").Value = 5500
$ws.Range("O3").Value = 6000
$ws.Range("P3").Value = 5750
$ws.Range("S3").Value = 1917

$ws.Range("D4").Value = 44172
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 160
$ws.Range("N4").Value = 5000
$ws.Range("O4").Value = 5500
$ws.Range("P4").Value = 5250
$ws.Range("S4").Value = 1750

$ws.Range("D5").Value = 44172
$ws.Range("L5").Value = "Tercera"
$ws.Range("M5").Value = 140
$ws.Range("N5").Value = 3500
$ws.Range("O5").Value = 4000
$ws.Range("P5").Value = 3750
$ws.Range("S5").Value = 1250

$ws.Range("D6").Value = 44895
$ws.Range("L6").Value = "Especial"
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 7000
$ws.Range("O6").Value = 8000
$ws.Range("P6").Value = 7500
$ws.Range("S6").Value = 2500

$ws.Range("D7").Value = 44895
$ws.Range("L7").Value = "Primera"
$ws.Range("M7").Value = 70
$ws.Range("N7").Value = 6000
$ws.Range("O7").Value = 7000
$ws.Range("P7").Value = 6500
$ws.Range("S7").Value = 2167

$ws.Range("D8").Value = 44895
$ws.Range("L8").Value = "Segunda"
$ws.Range("M8").Value = 72
$ws.Range("N8").Value = 5000
$ws.Range("O8").Value = 6000
$ws.Range("P8").Value = 5486
$ws.Range("S8").Value = 1829

$ws.Range("D9").Value = 44895
$ws.Range("L9").Value = "Tercera"
$ws.Range("M9").Value = 74
$ws.Range("N9").Value = 4000
$ws.Range("O9").Value = 5000
$ws.Range("P9").Value = 4500
$ws.Range("S9").Value = 1500

$ws.Range("D10").Value = 44334
$ws.Range("L10").Value = "Especial"
$ws.Range("M10").Value = 100
$ws.Range("N10").Value = 7000
$ws.Range("O10").Value = 8000
$ws.Range("P10").Value = 7500
$ws.Range("S10").Value = 2500

$ws.Range("D11").Value = 44334
$ws.Range("L11").Value = "Primera"
$ws.Range("M11").Value = 160
$ws.Range("N11").Value = 6000
$ws.Range("O11").Value = 7000
$ws.Range("P11").Value = 6500
$ws.Range("S11").Value = 2167

$ws.Range("D12").Value = 44334
$ws.Range("L12").Value = "Segunda"
$ws.Range("M12").Value = 120
$ws.Range("N12").Value = 6000
$ws.Range("O12").Value = 7000
$ws.Range("P12").Value = 6500
$ws.Range("S12").Value = 2167

$ws.Range("D13").Value = 44334
$ws.Range("L13").Value = "Tercera"
$ws.Range("M13").Value = 70
$ws.Range("N13").Value = 3500
$ws.Range("O13").Value = 4000
$ws.Range("P13").Value = 3750
$ws.Range("S13").Value = 1250

$ws.Range("D14").Value = 44811
$ws.Range("L14").Value = "Especial"
$ws.Range("M14").Value = 100
$ws.Range("N14").Value = 7000
$ws.Range("O14").Value = 8000
$ws.Range("P14").Value = 7500
$ws.Range("S14").Value = 2500

$ws.Range("D15").Value = 44811
$ws.Range("L15").Value = "Primera"
$ws.Range("M15").Value = 200
$ws.Range("N15").Value = 4000
$ws.Range("O15").Value = 5000
$ws.Range("P15").Value = 4500
$ws.Range("S15").Value = 1500

$ws.Range("D16").Value = 44811
$ws.Range("L16").Value = "Segunda"
$ws.Range("M16").Value = 200
$ws.Range("N16").Value = 3000
$ws.Range("O16").Value = 4000
$ws.Range("P16").Value = 3500
$ws.Range("S16").Value = 1167

$ws.Range("D17").Value = 44596
$ws.Range("L17").Value = "Especial"
$ws.Range("M17").Value = 100
$ws.Range("N17").Value = 8000
$ws.Range("O17").Value = 9000
$ws.Range("P17").Value = 8500
$ws.Range("S17").Value = 2833

$ws.Range("D18").Value = 44596
$ws.Range("L18").Value = "Primera"
$ws.Range("M18").Value = 130
$ws.Range("N18").Value = 6000
$ws.Range("O18").Value = 7000
$ws.Range("P18").Value = 6500
$ws.Range("S18").Value = 2167

$ws.Range("D19").Value = 44596
$ws.Range("L19").Value = "Segunda"
$ws.Range("M19").Value = 160
$ws.Range("N19").Value = 5000
$ws.Range("O19").Value = 6000
$ws.Range("P19").Value = 5500
$ws.Range("S19").Value = 1833

$ws.Range("D20").Value = 44596
$ws.Range("L20").Value = "Tercera"
$ws.Range("M20").Value = 100
$ws.Range("N20").Value = 4000
$ws.Range("O20").Value = 5000
$ws.Range("P20").Value = 4500
$ws.Range("S20").Value = 1500

$ws.Range("D21").Value = 44859
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 550
$ws.Range("N21").Value = 6000
$ws.Range("O21").Value = 7000
$ws.Range("P21").Value = 6545
$ws.Range("S21").Value = 2182

$ws.Range("D22").Value = 44859
$ws.Range("L22").Value = "Segunda"
$ws.Range("M22").Value = 500
$ws.Range("N22").Value = 5000
$ws.Range("O22").Value = 6000
$ws.Range("P22").Value = 5600
$ws.Range("S22").Value = 1867

$ws.Range("D23").Value = 44859
$ws.Range("L23").Value = "Tercera"
$ws.Range("M23").Value = 350
$ws.Range("N23").Value = 4000
$ws.Range("O23").Value = 5000
$ws.Range("P23").Value = 4857
$ws.Range("S23").Value = 1619

$ws.Range("D24").Value = 44242
$ws.Range("L24").Value = "Especial"
$ws.Range("M24").Value = 50
$ws.Range("N24").Value = 7000
$ws.Range("O24").Value = 8000
$ws.Range("P24").Value = 7500
$ws.Range("S24").Value = 2500

$ws.Range("D25").Value = 44242
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 90
$ws.Range("N25").Value = 6000
$ws.Range("O25").Value = 7000
$ws.Range("P25").Value = 6500
$ws.Range("S25").Value = 2167

$ws.Range("D26").Value = 44242
$ws.Range("L26").Value = "Segunda"
$ws.Range("M26").Value = 100
$ws.Range("N26").Value = 4000
$ws.Range("O26").Value = 5000
$ws.Range("P26").Value = 4500
$ws.Range("S26").Value = 1500

$ws.Range("D27").Value = 44389
$ws.Range("L27").Value = "Especial"
$ws.Range("M27").Value = 100
$ws.Range("N27").Value = 7500
$ws.Range("O27").Value = 8000
$ws.Range("P27").Value = 7750
$ws.Range("S27").Value = 2583

$ws.Range("D28").Value = 44389
$ws.Range("L28").Value = "Primera"
$ws.Range("M28").Value = 160
$ws.Range("N28").Value = 6000
$ws.Range("O28").Value = 7000
$ws.Range("P28").Value = 6500
$ws.Range("S28").Value = 2167

$ws.Range("D29").Value = 44389
$ws.Range("L29").Value = "Segunda"
$ws.Range("M29").Value = 200
$ws.Range("N29").Value = 5500
$ws.Range("O29").Value = 6000
$ws.Range("P29").Value = 5750
$ws.Range("S29").Value = 1917

$ws.Range("D30").Value = 44708
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 50
$ws.Range("N30").Value = 6000
$ws.Range("O30").Value = 7000
$ws.Range("P30").Value = 6500
$ws.Range("S30").Value = 2167

$ws.Range("D31").Value = 44708
$ws.Range("L31").Value = "Segunda"
$ws.Range("M31").Value = 60
$ws.Range("N31").Value = 4000
$ws.Range("O31").Value = 5000
$ws.Range("P31").Value = 4500
$ws.Range("S31").Value = 1500

$ws.Range("D32").Value = 44708
$ws.Range("L32").Value = "Tercera"
$ws.Range("M32").Value = 50
$ws.Range("N32").Value = 3000
$ws.Range("O32").Value = 4000
$ws.Range("P32").Value = 3500
$ws.Range("S32").Value = 1167

$ws.Range("D33").Value = 44894
$ws.Range("L33").Value = "Especial"
$ws.Range("M33").Value = 60
$ws.Range("N33").Value = 7000
$ws.Range("O33").Value = 8000
$ws.Range("P33").Value = 7500
$ws.Range("S33").Value = 2500

$ws.Range("D34").Value = 44894
$ws.Range("L34").Value = "Primera"
$ws.Range("M34").Value = 70
$ws.Range("N34").Value = 6000
$ws.Range("O34").Value = 7000
$ws.Range("P34").Value = 6500
$ws.Range("S34").Value = 2167

$ws.Range("D35").Value = 44894
$ws.Range("L35").Value = "Segunda"
$ws.Range("M35").Value = 72
$ws.Range("N35").Value = 5000
$ws.Range("O35").Value = 6000
$ws.Range("P35").Value = 5486
$ws.Range("S35").Value = 1829

$ws.Range("D36").Value = 44894
$ws.Range("L36").Value = "Tercera"
$ws.Range("M36").Value = 74
$ws.Range("N36").Value = 4000
$ws.Range("O36").Value = 5000
$ws.Range("P36").Value = 4500
$ws.Range("S36").Value = 1500

$ws.Range("D37").Value = 44200
$ws.Range("L37").Value = "Especial"
$ws.Range("M37").Value = 50
$ws.Range("N37").Value = 4500
$ws.Range("O37").Value = 5000
$ws.Range("P37").Value = 4750
$ws.Range("S37").Value = 1583

$ws.Range("D38").Value = 44200
$ws.Range("L38").Value = "Primera"
$ws.Range("M38").Value = 80
$ws.Range("N38").Value = 3500
$ws.Range("O38").Value = 4000
$ws.Range("P38").Value = 3750
$ws.Range("S38").Value = 1250

$ws.Range("D39").Value = 44200
$ws.Range("L39").Value = "Segunda"
$ws.Range("M39").Value = 120
$ws.Range("N39").Value = 2500
$ws.Range("O39").Value = 3000
$ws.Range("P39").Value = 2750
$ws.Range("S39").Value = 917

$ws.Range("D40").Value = 44249
$ws.Range("L40").Value = "Especial"
$ws.Range("M40").Value = 200
$ws.Range("N40").Value = 6000
$ws.Range("O40").Value = 7000
$ws.Range("P40").Value = 6500
$ws.Range("S40").Value = 2167

$ws.Range("D41").Value = 44249
$ws.Range("L41").Value = "Primera"
$ws.Range("M41").Value = 160
$ws.Range("N41").Value = 4500
$ws.Range("O41").Value = 5000
$ws.Range("P41").Value = 4750
$ws.Range("S41").Value = 1583

$ws.Range("D42").Value = 44832
$ws.Range("L42").Value = "Especial"
$ws.Range("M42").Value = 100
$ws.Range("N42").Value = 6000
$ws.Range("O42").Value = 7000
$ws.Range("P42").Value = 6500
$ws.Range("S42").Value = 2167

$ws.Range("D43").Value = 44832
$ws.Range("L43").Value = "Primera"
$ws.Range("M43").Value = 120
$ws.Range("N43").Value = 5000
$ws.Range("O43").Value = 6000
$ws.Range("P43").Value = 5500
$ws.Range("S43").Value = 1833

$ws.Range("D44").Value = 44832
$ws.Range("L44").Value = "Segunda"
$ws.Range("M44").Value = 140
$ws.Range("N44").Value = 4000
$ws.Range("O44").Value = 5000
$ws.Range("P44").Value = 4500
$ws.Range("S44").Value = 1500

$ws.Range("D45").Value = 44855
$ws.Range("L45").Value = "Especial"
$ws.Range("M45").Value = 160
$ws.Range("N45").Value = 6500
$ws.Range("O45").Value = 7000
$ws.Range("P45").Value = 6750
$ws.Range("S45").Value = 2250

$ws.Range("D46").Value = 44855
$ws.Range("L46").Value = "Primera"
$ws.Range("M46").Value = 160
$ws.Range("N46").Value = 5500
$ws.Range("O46").Value = 6000
$ws.Range("P46").Value = 5750
$ws.Range("S46").Value = 1917

$ws.Range("D47").Value = 44855
$ws.Range("L47").Value = "Segunda"
$ws.Range("M47").Value = 200
$ws.Range("N47").Value = 4500
$ws.Range("O47").Value = 5000
$ws.Range("P47").Value = 4750
$ws.Range("S47").Value = 1583

$ws.Range("D48").Value = 44322
$ws.Range("L48").Value = "Especial"
$ws.Range("M48").Value = 200
$ws.Range("N48").Value = 7000
$ws.Range("O48").Value = 7500
$ws.Range("P48").Value = 7250
$ws.Range("S48").Value = 2417

$ws.Range("D49").Value = 44322
$ws.Range("L49").Value = "Primera"
$ws.Range("M49").Value = 160
$ws.Range("N49").Value = 6000
$ws.Range("O49").Value = 6500
$ws.Range("P49").Value = 6250
$ws.Range("S49").Value = 2083

$ws.Range("D50").Value = 44322
$ws.Range("L50").Value = "Segunda"
$ws.Range("M50").Value = 100
$ws.Range("N50").Value = 5000
$ws.Range("O50").Value = 5500
$ws.Range("P50").Value = 5250
$ws.Range("S50").Value = 1750

$ws.Range("D51").Value = 44351
$ws.Range("L51").Value = "Especial"
$ws.Range("M51").Value = 160
$ws.Range("N51").Value = 7500
$ws.Range("O51").Value = 8000
$ws.Range("P51").Value = 7750
$ws.Range("S51").Value = 2583

$ws.Range("D52").Value = 44351
$ws.Range("L52").Value = "Primera"
$ws.Range("M52").Value = 100
$ws.Range("N52").Value = 6000
$ws.Range("O52").Value = 6500
$ws.Range("P52").Value = 6250
$ws.Range("S52").Value = 2083

$ws.Range("D53").Value = 44351
$ws.Range("L53").Value = "Segunda"
$ws.Range("M53").Value = 200
$ws.Range("N53").Value = 4500
$ws.Range("O53").Value = 5000
$ws.Range("P53").Value = 4750
$ws.Range("S53").Value = 1583
